$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (blad1 -> 06020)
$ws.Name = "06020"

# Clear the rotated-text style previously applied to row 1 (A1:G1), so the
# workbook no longer needs the extra cellXfs entry with textRotation=45.
$ws.Range("A1:G1").Style = "Normal"
$ws.Range("A1:G1").Clear()

# Row 1: header labels
$ws.Range("A1").Value = "omschrijving"
$ws.Range("B1").Value = "inhoud"
$ws.Range("C1").Value = "weergave"
$ws.Range("D1").Value = "uitlijnen"
$ws.Range("E1").Value = "regel verwijderen"
$ws.Range("F1").Value = "regel template"
$ws.Range("G1").Value = "PTEST"

# Row 2
$ws.Range("A2").Value = "Eigen risico"
$ws.Range("B2").Value = "€ 10043"
$ws.Range("F2").Value = "06 Eigen risico               € 10043"
$ws.Range("G2").Value = "x"

# Row 3
$ws.Range("B3").Value = "10043"
$ws.Range("C3").Value = "Getal inclusief decimalen"
$ws.Range("D3").Value = "Links"
$ws.Range("E3").Value = "verwijderen"

# Move the selection to G2 (matches the saved cursor position in the file)
$ws.Range("G2").Select()
